$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B18").Value = 18
$ws.Range("B18").Select()
